$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update F column (想去人数 / "people who want to go")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1035
$ws1.Range("F4").Value = 173
$ws1.Range("F5").Value = 2817
$ws1.Range("F7").Value = 230
$ws1.Range("F9").Value = 124
$ws1.Range("F10").Value = 77
$ws1.Range("F11").Value = 92
$ws1.Range("F12").Value = 2646
$ws1.Range("F13").Value = 834

# Sheet "全部类型" (All types) - same underlying events, offset by one row
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1035
$ws4.Range("F5").Value = 173
$ws4.Range("F6").Value = 2817
$ws4.Range("F8").Value = 230
$ws4.Range("F11").Value = 124
$ws4.Range("F12").Value = 77
$ws4.Range("F13").Value = 92
$ws4.Range("F14").Value = 2646
$ws4.Range("F15").Value = 834
